$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking columns A, B, E so Excel keeps them as text
$ws.Range("A52:B92").NumberFormat = "@"
$ws.Range("E52:E92").NumberFormat = "@"

$ws.Range('A52').Value = '0'
$ws.Range('B52').Value = '1222440'
$ws.Range('C52').Value = '"Sneaky Ninja Panda"'
$ws.Range('D52').Value = '二馆'
$ws.Range('E52').Value = '0'

$ws.Range('A53').Value = '0'
$ws.Range('B53').Value = '3391765'
$ws.Range('C53').Value = '马er'
$ws.Range('D53').Value = '二馆'
$ws.Range('E53').Value = '0'

$ws.Range('A54').Value = '69505'
$ws.Range('B54').Value = '9718882'
$ws.Range('C54').Value = '小霸王2021'
$ws.Range('D54').Value = '二馆'
$ws.Range('E54').Value = '2517'

$ws.Range('A55').Value = '44146'
$ws.Range('B55').Value = '11645391'
$ws.Range('C55').Value = '"omar omar"'
$ws.Range('D55').Value = '二馆'
$ws.Range('E55').Value = '3936'

$ws.Range('A56').Value = '0'
$ws.Range('B56').Value = '15436348'
$ws.Range('C56').Value = 'Lucas'
$ws.Range('D56').Value = '二馆'
$ws.Range('E56').Value = '1525'

$ws.Range('A57').Value = '0'
$ws.Range('B57').Value = '20372140'
$ws.Range('C57').Value = '人山即是仙'
$ws.Range('D57').Value = '二馆'
$ws.Range('E57').Value = '0'

$ws.Range('A58').Value = '0'
$ws.Range('B58').Value = '38994054'
$ws.Range('C58').Value = 'chengnan'
$ws.Range('D58').Value = '二馆'
$ws.Range('E58').Value = '0'

$ws.Range('A59').Value = '61497'
$ws.Range('B59').Value = '41837764'
$ws.Range('C59').Value = '好风光会长'
$ws.Range('D59').Value = '二馆'
$ws.Range('E59').Value = '2689'

$ws.Range('A60').Value = '0'
$ws.Range('B60').Value = '43281368'
$ws.Range('C60').Value = 'xhs2763'
$ws.Range('D60').Value = '二馆'
$ws.Range('E60').Value = '0'

$ws.Range('A61').Value = '0'
$ws.Range('B61').Value = '44378757'
$ws.Range('C61').Value = '"NᵉᵗʰᵉʳDʳⁱᶠᵗᵉʳ ㊥"'
$ws.Range('D61').Value = '二馆'
$ws.Range('E61').Value = '0'

$ws.Range('A62').Value = '43183'
$ws.Range('B62').Value = '47430231'
$ws.Range('C62').Value = 'Kentantrino'
$ws.Range('D62').Value = '二馆'
$ws.Range('E62').Value = '4002'

$ws.Range('A63').Value = '50821'
$ws.Range('B63').Value = '48738257'
$ws.Range('C63').Value = '死亡洲际跳蛋'
$ws.Range('D63').Value = '二馆'
$ws.Range('E63').Value = '3209'

$ws.Range('A64').Value = '0'
$ws.Range('B64').Value = '49000199'
$ws.Range('C64').Value = 'SlipperyForester5672'
$ws.Range('D64').Value = '二馆'
$ws.Range('E64').Value = '1300'

$ws.Range('A65').Value = '0'
$ws.Range('B65').Value = '54941706'
$ws.Range('C65').Value = 'AlexMenjivar20'
$ws.Range('D65').Value = '二馆'
$ws.Range('E65').Value = '1470'

$ws.Range('A66').Value = '60472'
$ws.Range('B66').Value = '55499394'
$ws.Range('C66').Value = 'Player-55499394'
$ws.Range('D66').Value = '二馆'
$ws.Range('E66').Value = '2722'

$ws.Range('A67').Value = '0'
$ws.Range('B67').Value = '55810157'
$ws.Range('C67').Value = 'Beard'
$ws.Range('D67').Value = '二馆'
$ws.Range('E67').Value = '0'

$ws.Range('A68').Value = '0'
$ws.Range('B68').Value = '56700848'
$ws.Range('C68').Value = '工口漫画老师'
$ws.Range('D68').Value = '二馆'
$ws.Range('E68').Value = '0'

$ws.Range('A69').Value = '0'
$ws.Range('B69').Value = '57219176'
$ws.Range('C69').Value = '青莲道人'
$ws.Range('D69').Value = '二馆'
$ws.Range('E69').Value = '1494'

$ws.Range('A70').Value = '0'
$ws.Range('B70').Value = '57556179'
$ws.Range('C70').Value = '特战新生代英雄'
$ws.Range('D70').Value = '二馆'
$ws.Range('E70').Value = '0'

$ws.Range('A71').Value = '0'
$ws.Range('B71').Value = '58340439'
$ws.Range('C71').Value = '70qilin'
$ws.Range('D71').Value = '二馆'
$ws.Range('E71').Value = '0'

$ws.Range('A72').Value = '0'
$ws.Range('B72').Value = '58615925'
$ws.Range('C72').Value = '齐天的大圣'
$ws.Range('D72').Value = '二馆'
$ws.Range('E72').Value = '0'

$ws.Range('A73').Value = '0'
$ws.Range('B73').Value = '58641574'
$ws.Range('C73').Value = 'Player-58641574鱼'
$ws.Range('D73').Value = '二馆'
$ws.Range('E73').Value = '0'

$ws.Range('A74').Value = '0'
$ws.Range('B74').Value = '58743790'
$ws.Range('C74').Value = 'Ma'
$ws.Range('D74').Value = '二馆'
$ws.Range('E74').Value = '0'

$ws.Range('A75').Value = '98788'
$ws.Range('B75').Value = '8666978'
$ws.Range('C75').Value = 'FierceRocket'
$ws.Range('D75').Value = '三馆'
$ws.Range('E75').Value = '1633'

$ws.Range('A76').Value = '0'
$ws.Range('B76').Value = '15695258'
$ws.Range('C76').Value = 'Player-15695258'
$ws.Range('D76').Value = '三馆'
$ws.Range('E76').Value = '1000'

$ws.Range('A77').Value = '0'
$ws.Range('B77').Value = '29355299'
$ws.Range('C77').Value = 'Player-29355299'
$ws.Range('D77').Value = '三馆'
$ws.Range('E77').Value = '0'

$ws.Range('A78').Value = '0'
$ws.Range('B78').Value = '41231396'
$ws.Range('C78').Value = 'ollsthebro'
$ws.Range('D78').Value = '三馆'
$ws.Range('E78').Value = '0'

$ws.Range('A79').Value = '0'
$ws.Range('B79').Value = '47622456'
$ws.Range('C79').Value = '伊恩'
$ws.Range('D79').Value = '三馆'
$ws.Range('E79').Value = '0'

$ws.Range('A80').Value = '0'
$ws.Range('B80').Value = '49553719'
$ws.Range('C80').Value = '"Oreo Captain Sir"'
$ws.Range('D80').Value = '三馆'
$ws.Range('E80').Value = '0'

$ws.Range('A81').Value = '0'
$ws.Range('B81').Value = '55745105'
$ws.Range('C81').Value = 'eldeniz'
$ws.Range('D81').Value = '三馆'
$ws.Range('E81').Value = '0'

$ws.Range('A82').Value = '47528'
$ws.Range('B82').Value = '56241637'
$ws.Range('C82').Value = 'Player-14day'
$ws.Range('D82').Value = '三馆'
$ws.Range('E82').Value = '3509'

$ws.Range('A83').Value = '0'
$ws.Range('B83').Value = '58174442'
$ws.Range('C83').Value = 'Player-58174442'
$ws.Range('D83').Value = '三馆'
$ws.Range('E83').Value = '1020'

$ws.Range('A84').Value = '0'
$ws.Range('B84').Value = '58572199'
$ws.Range('C84').Value = '你干嘛～哎呦～'
$ws.Range('D84').Value = '三馆'
$ws.Range('E84').Value = '0'

$ws.Range('A85').Value = '0'
$ws.Range('B85').Value = '58671339'
$ws.Range('C85').Value = '"quang pro"'
$ws.Range('D85').Value = '三馆'
$ws.Range('E85').Value = '0'

$ws.Range('A86').Value = '0'
$ws.Range('B86').Value = '58766144'
$ws.Range('C86').Value = 'EquablePrecedence38'
$ws.Range('D86').Value = '三馆'
$ws.Range('E86').Value = '0'

$ws.Range('A87').Value = '0'
$ws.Range('B87').Value = '58910668'
$ws.Range('C87').Value = 'BrittleAuthor33'
$ws.Range('D87').Value = '三馆'
$ws.Range('E87').Value = '0'

$ws.Range('A88').Value = '0'
$ws.Range('B88').Value = '59081265'
$ws.Range('C88').Value = '爬楼梯'
$ws.Range('D88').Value = '三馆'
$ws.Range('E88').Value = '0'

$ws.Range('A89').Value = '0'
$ws.Range('B89').Value = '59082827'
$ws.Range('C89').Value = 'Player-59082827'
$ws.Range('D89').Value = '三馆'
$ws.Range('E89').Value = '0'

$ws.Range('A90').Value = '0'
$ws.Range('B90').Value = '59106471'
$ws.Range('C90').Value = 'anime'
$ws.Range('D90').Value = '三馆'
$ws.Range('E90').Value = '0'

$ws.Range('A91').Value = '0'
$ws.Range('B91').Value = '59112086'
$ws.Range('C91').Value = 'sigma'
$ws.Range('D91').Value = '三馆'
$ws.Range('E91').Value = '0'

$ws.Range('A92').Value = '91292'
$ws.Range('B92').Value = '59304163'
$ws.Range('C92').Value = 'Hong'
$ws.Range('D92').Value = '三馆'
$ws.Range('E92').Value = '1909'

# Remove the now-obsolete trailing rows (old rows 93-115)
$ws.Range("A93:E115").EntireRow.Delete() | Out-Null